$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: quote retrieved -> ExpectedRate matches ActualRate, Result = PASS
$ws.Range("E2").Value = "'$19.04"
$ws.Range("E2").ClearFormats()

$ws.Range("F2").Value = "PASS"

# Row 30: new quote retrieved -> ActualRate differs from before, Result = FAIL
$ws.Range("E30").Value = "'$473.23"
$ws.Range("E30").ClearFormats()

$ws.Range("F30").Value = "FAIL"
